$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 211.4614666666667
$ws.Range("H2").Value = 634.3844
$ws.Range("I2").Value = 0.2421062275331183
$ws.Range("J2").Value = 0.2421062275331183
$ws.Range("O2").Value = 0.9418062875790357
$ws.Range("P2").Value = 0.9418062875790357
$ws.Range("Q2").Value = 43.65586735755556
$ws.Range("R2").Value = 392.902806218
$ws.Range("S2").Value = 0.2280171673527315
$ws.Range("T2").Value = 0.2280171673527315

# Row 3
$ws.Range("G3").Value = 211.4614666666667
$ws.Range("H3").Value = 634.3844
$ws.Range("I3").Value = 0.2421062275331183
$ws.Range("J3").Value = 0.2421062275331183
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01275633333333333
$ws.Range("N3").Value = 0.038269
$ws.Range("O3").Value = 0.05819371242096427
$ws.Range("P3").Value = 0.05819371242096427
$ws.Range("Q3").Value = 2.697472955955555
$ws.Range("R3").Value = 24.2772566036
$ws.Range("S3").Value = 0.01408906018038683
$ws.Range("T3").Value = 0.01408906018038683

# Row 4
$ws.Range("I4").Value = 0.08842543241393927
$ws.Range("J4").Value = 0.08842543241393927
$ws.Range("O4").Value = 0.9418062875790357
$ws.Range("P4").Value = 0.9418062875790357
$ws.Range("S4").Value = 0.08327962822934308
$ws.Range("T4").Value = 0.08327962822934308

# Row 5
$ws.Range("I5").Value = 0.08842543241393927
$ws.Range("J5").Value = 0.08842543241393927
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01275633333333333
$ws.Range("N5").Value = 0.038269
$ws.Range("O5").Value = 0.05819371242096427
$ws.Range("P5").Value = 0.05819371242096427
$ws.Range("Q5").Value = 0.9852089100956667
$ws.Range("R5").Value = 8.866880190861
$ws.Range("S5").Value = 0.005145804184596194
$ws.Range("T5").Value = 0.005145804184596194

# Row 6
$ws.Range("G6").Value = 174.3107043333333
$ws.Range("H6").Value = 522.932113
$ws.Range("I6").Value = 0.199571617988009
$ws.Range("J6").Value = 0.199571617988009
$ws.Range("O6").Value = 0.9418062875790357
$ws.Range("P6").Value = 0.9418062875790357
$ws.Range("Q6").Value = 35.98615439177611
$ws.Range("R6").Value = 323.875389525985
$ws.Range("S6").Value = 0.1879578046434283
$ws.Range("T6").Value = 0.1879578046434283

# Row 7
$ws.Range("G7").Value = 174.3107043333333
$ws.Range("H7").Value = 522.932113
$ws.Range("I7").Value = 0.199571617988009
$ws.Range("J7").Value = 0.199571617988009
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01275633333333333
$ws.Range("N7").Value = 0.038269
$ws.Range("O7").Value = 0.05819371242096427
$ws.Range("P7").Value = 0.05819371242096427
$ws.Range("Q7").Value = 2.223565448044111
$ws.Range("R7").Value = 20.012089032397
$ws.Range("S7").Value = 0.01161381334458074
$ws.Range("T7").Value = 0.01161381334458074

# Row 8
$ws.Range("G8").Value = 28.53474833333333
$ws.Range("H8").Value = 85.60424499999999
$ws.Range("I8").Value = 0.03266997236655063
$ws.Range("J8").Value = 0.03266997236655063
$ws.Range("O8").Value = 0.9418062875790357
$ws.Range("P8").Value = 0.9418062875790357
$ws.Range("Q8").Value = 5.890951235502778
$ws.Range("R8").Value = 53.018561119525
$ws.Range("S8").Value = 0.03076878538985073
$ws.Range("T8").Value = 0.03076878538985073

# Row 9
$ws.Range("G9").Value = 28.53474833333333
$ws.Range("H9").Value = 85.60424499999999
$ws.Range("I9").Value = 0.03266997236655063
$ws.Range("J9").Value = 0.03266997236655063
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01275633333333333
$ws.Range("N9").Value = 0.038269
$ws.Range("O9").Value = 0.05819371242096427
$ws.Range("P9").Value = 0.05819371242096427
$ws.Range("Q9").Value = 0.3639987613227777
$ws.Range("R9").Value = 3.275988851904999
$ws.Range("S9").Value = 0.001901186976699897
$ws.Range("T9").Value = 0.001901186976699897

# Row 10
$ws.Range("G10").Value = 230.32901
$ws.Range("H10").Value = 690.98703
$ws.Range("I10").Value = 0.263708034289011
$ws.Range("J10").Value = 0.263708034289011
$ws.Range("O10").Value = 0.9418062875790357
$ws.Range("P10").Value = 0.9418062875790357
$ws.Range("Q10").Value = 47.55104023281667
$ws.Range("R10").Value = 427.95936209535
$ws.Range("S10").Value = 0.2483618847784985
$ws.Range("T10").Value = 0.2483618847784985

# Row 11
$ws.Range("G11").Value = 230.32901
$ws.Range("H11").Value = 690.98703
$ws.Range("I11").Value = 0.263708034289011
$ws.Range("J11").Value = 0.263708034289011
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.01275633333333333
$ws.Range("N11").Value = 0.038269
$ws.Range("O11").Value = 0.05819371242096427
$ws.Range("P11").Value = 0.05819371242096427
$ws.Range("Q11").Value = 2.938153627896667
$ws.Range("R11").Value = 26.44338265107
$ws.Range("S11").Value = 0.01534614951051249
$ws.Range("T11").Value = 0.01534614951051249

# Row 12
$ws.Range("G12").Value = 151.5554656666667
$ws.Range("H12").Value = 454.666397
$ws.Range("I12").Value = 0.1735187154093718
$ws.Range("J12").Value = 0.1735187154093718
$ws.Range("O12").Value = 0.9418062875790357
$ws.Range("P12").Value = 0.9418062875790357
$ws.Range("Q12").Value = 31.28837329444056
$ws.Range("R12").Value = 281.595359649965
$ws.Range("S12").Value = 0.1634210171851837
$ws.Range("T12").Value = 0.1634210171851837

# Row 13
$ws.Range("G13").Value = 151.5554656666667
$ws.Range("H13").Value = 454.666397
$ws.Range("I13").Value = 0.1735187154093718
$ws.Range("J13").Value = 0.1735187154093718
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.01275633333333333
$ws.Range("N13").Value = 0.038269
$ws.Range("O13").Value = 0.05819371242096427
$ws.Range("P13").Value = 0.05819371242096427
$ws.Range("Q13").Value = 1.933292038532555
$ws.Range("R13").Value = 17.399628346793
$ws.Range("S13").Value = 0.01009769822418812
$ws.Range("T13").Value = 0.01009769822418812

